$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the shared string text typos (remove trailing "2")
$ws.Range("D1").Value = "Spanish"
$ws.Range("E1").Value = "Portuguese"

# Clear the leftover D2 cell entirely (value + formatting)
$ws.Range("D2").ClearContents()
$ws.Range("D2").ClearFormats()

# Update the current selection on the sheet
$ws.Range("D2").Select()

# Adjust the default column width slightly
$ws.Columns.ColumnWidth = 11.55078125
